# "update bat and fish" — fill in LEFT/TOP/RIGHT/BOTTOM sprite-sheet coordinates
# for door/gate items 10230-10232 on sheet1 (Trang_tinh1), label the first one,
# and add two corresponding composite-sprite rows (320/321) on sheet2 (Trang_tinh2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1: Trang_tinh1 ------------------------------------------------
# Row 31 (ID 10230): LEFT/TOP/RIGHT/BOTTOM + status label "cua scene 2.1"
$ws1.Range("H31").Value = 700
$ws1.Range("I31").Value = 130
$ws1.Range("J31").Value = 716
$ws1.Range("K31").Value = 195
$ws1.Range("L31").Value = "cửa scene 2.1"

# Row 32 (ID 10231): LEFT/TOP/RIGHT/BOTTOM
$ws1.Range("H32").Value = 768
$ws1.Range("I32").Value = 130
$ws1.Range("J32").Value = 800
$ws1.Range("K32").Value = 195

# Row 33 (ID 10232): LEFT/TOP/RIGHT/BOTTOM
$ws1.Range("H33").Value = 800
$ws1.Range("I33").Value = 130
$ws1.Range("J33").Value = 900
$ws1.Range("K33").Value = 195

# ---- Sheet2: Trang_tinh2 -------------------------------------------------
# Row 22 (ID 320): sprite refs 10230/10231/10232/10231 + status "cua 2"
$ws2.Range("A22").Value = 320
$ws2.Range("B22").Value = 10230
$ws2.Range("C22").Value = 10231
$ws2.Range("D22").Value = 10232
$ws2.Range("E22").Value = 10231
$ws2.Range("F22").Value = "cửa 2"

# Row 23 (ID 321): sprite refs 10230/-1/-1/-1
$ws2.Range("A23").Value = 321
$ws2.Range("B23").Value = 10230
$ws2.Range("C23").Value = -1
$ws2.Range("D23").Value = -1
$ws2.Range("E23").Value = -1

# ---- View state: match the final selections from the edit -------------
# Sheet1: scroll the frozen pane and leave G31:K33 selected.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("G31:K33").Select() | Out-Null

# Sheet2 ends up the active tab/selection, matching the source file.
$ws2.Activate() | Out-Null
$ws2.Range("A22:E23").Select() | Out-Null
